$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.595.81'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '1.924.42'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.00'
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4815'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4061'
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08224'
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.011'
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.72'
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("D12").Value = '1.932.15'
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.077'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.262'
$ws.Range("E14").Value = '  +1.61%  '
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06874'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("E20").Value = '  +0.40%  '
$ws.Range("D21").Value = '29.592.62'
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.682'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.96'
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.184'
$ws.Range("D25").Value = '2.160.97'
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.95'
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.469'
$ws.Range("E27").Value = '  +0.79%  '
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.092'
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.016'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09631'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.620'
$ws.Range("E33").Value = '  +2.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.576'
$ws.Range("E34").Value = '  +0.25%  '
$ws.Range("E35").Value = '  -1.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06336'
$ws.Range("E36").Value = '  +3.81%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02292'
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5960'
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("E40").Value = '  +0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.881'
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1853'
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.288'
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.41'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07505'
$ws.Range("E46").Value = '  -2.26%  '
$ws.Range("E47").Value = '  -0.41%  '
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.58'
$ws.Range("E49").Value = '  +3.41%  '
$ws.Range("E50").Value = '  +3.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.19'
$ws.Range("E51").Value = '  -0.66%  '
